# Update the "writing" progress log with the latest day's entry (11/11/2020)
# and refresh the dependent Table1 + dashboard chart, matching the author's
# "update plots; create slides" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("writing")

# --- 1. Grow Table1 by one row (this extends both the table ref and the
#        worksheet dimension automatically) ---------------------------------
$lo = $ws.ListObjects.Item("Table1")
$lo.ListRows.Add() | Out-Null

# Copy the date formatting (m/d/yyyy) from the row above down into the new
# row's date cell before putting the value in, so it keeps the same number
# format / style as the rest of column A.
$ws.Range("A13").Copy()
$ws.Range("A14").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- 2. Fill in the new day's data (row 14) --------------------------------
$ws.Range("A14").Value = 44146
$ws.Range("B14").Value = 239
$ws.Range("C14").Value = 87
$ws.Range("D14").Value = 547
$ws.Range("E14").Value = 7413
$ws.Range("F14").Value = 401
$ws.Range("G14").Value = 120
$ws.Range("H14").Value = 117
$ws.Range("I14").Value = 186
$ws.Range("J14").Formula = "=SUM(B14:I14)"
$ws.Range("K14").Value = 2671

# --- 3. Re-enter the "Total" column formula across J8:J13 as one fill so
#        Excel stores it as a shared formula group (matches rows 2:7) ------
$ws.Range("J8:J13").Formula = "=SUM(B8:I8)"

# --- 4. Point the dashboard chart's two series at the new A2:A14/J2:J14/
#        K2:K14 ranges ------------------------------------------------------
$dash = $wb.Worksheets.Item("dashboard")
$chartObj = $dash.ChartObjects().Item(1)
$chart = $chartObj.Chart

$serDaily = $chart.SeriesCollection().Item(1)   # "Daily" bars -> writing!$K$2:$K$13
$serDaily.Values = "=writing!`$K`$2:`$K`$14"
$serDaily.XValues = "=writing!`$A`$2:`$A`$14"

$serTotal = $chart.SeriesCollection().Item(2)   # "Total" line -> writing!$J$2:$J$13
$serTotal.Values = "=writing!`$J`$2:`$J`$14"
$serTotal.XValues = "=writing!`$A`$2:`$A`$14"

# --- 5. Match the author's final UI state: selection sitting on K14 of the
#        "writing" sheet, but "dashboard" left as the active/visible tab ---
$ws.Range("K14").Select()
$dash.Activate()
